$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11: new log entry ---
$ws.Range("A11").Value = 46072
$ws.Range("A11").NumberFormat = "d-mmm"
$ws.Range("A11").WrapText = $true

$ws.Range("B11").Value = "Bit Manipulation Questions"
$ws.Range("C11").Value = "None"
$ws.Range("D11").Value = "100 mins"

# --- Row 12: new log entry ---
$ws.Range("A12").Value = 46073
$ws.Range("A12").NumberFormat = "d-mmm"
$ws.Range("A12").WrapText = $true

$ws.Range("B12").Value = "Compled bit-manipulation and started OOP concepts"
$ws.Range("C12").Value = "count bits - 338"
$ws.Range("D12").Value = "150-180 mins"

# Row heights grow to fit the wrapped text in the new rows
$ws.Rows.Item(11).RowHeight = 28.8
$ws.Rows.Item(12).RowHeight = 57.6

# Move the active selection to D13, where the cursor ended up after logging
[void]$ws.Range("D13").Select()
